$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  'D2' = 44676
  'M2' = 55
  'N2' = 28000
  'O2' = 30000
  'P2' = 28909
  'S2' = 1606
  'D3' = 44629
  'D4' = 44424
  'R4' = 'Región de Arica y Parinacota'
  'D5' = 44405
  'M5' = 10
  'D6' = 44392
  'M6' = 20
  'N6' = 35000
  'O6' = 35000
  'P6' = 35000
  'Q6' = '$/caja 18 kilos'
  'S6' = 1944
  'T6' = 18
  'D7' = 44369
  'M7' = 5
  'N7' = 35000
  'O7' = 35000
  'P7' = 35000
  'S7' = 1944
  'D8' = 44719
  'M8' = 25
  'N8' = 34000
  'O8' = 34000
  'P8' = 34000
  'S8' = 1889
  'D9' = 44721
  'M9' = 5
  'N9' = 35000
  'O9' = 35000
  'P9' = 35000
  'R9' = 'Perú'
  'S9' = 1944
  'D10' = 44364
  'M10' = 90
  'N10' = 1700
  'O10' = 1700
  'P10' = 1700
  'Q10' = '$/kilo'
  'S10' = 1700
  'T10' = 1
  'D11' = 44264
  'M11' = 20
  'N11' = 40000
  'O11' = 40000
  'P11' = 40000
  'S11' = 2222
  'D12' = 44669
  'M12' = 40
  'N12' = 32000
  'O12' = 32000
  'P12' = 32000
  'R12' = 'Región de Arica y Parinacota'
  'S12' = 1778
  'D13' = 44664
  'M13' = 15
  'N13' = 30000
  'O13' = 30000
  'P13' = 30000
  'S13' = 1667
  'D14' = 44363
  'M14' = 144
  'N14' = 1700
  'O14' = 1700
  'P14' = 1700
  'Q14' = '$/kilo'
  'S14' = 1700
  'T14' = 1
  'D15' = 44431
  'O15' = 35000
  'P15' = 35000
  'S15' = 1944
  'D16' = 44671
  'M16' = 20
  'N16' = 32000
  'O16' = 32000
  'P16' = 32000
  'S16' = 1778
  'D17' = 44634
  'M17' = 30
  'N17' = 45000
  'O17' = 45000
  'P17' = 45000
  'S17' = 2500
  'D18' = 44645
  'M18' = 5
  'N18' = 30000
  'O18' = 30000
  'P18' = 30000
  'S18' = 1667
  'D19' = 44438
  'M19' = 25
  'D20' = 44704
  'M20' = 25
  'D21' = 44434
  'M21' = 40
  'N21' = 35000
  'O21' = 35000
  'P21' = 35000
  'S21' = 1944
  'D22' = 44448
  'M22' = 50
  'N22' = 38000
  'O22' = 38000
  'P22' = 38000
  'R22' = 'Región de Arica y Parinacota'
  'S22' = 2111
  'D23' = 44279
  'M23' = 30
  'N23' = 35000
  'O23' = 36000
  'P23' = 35667
  'S23' = 1982
  'D24' = 44679
  'M24' = 35
  'N24' = 34000
  'O24' = 34000
  'P24' = 34000
  'R24' = 'Perú'
  'S24' = 1889
  'D25' = 44679
  'M25' = 55
  'N25' = 28000
  'O25' = 28000
  'P25' = 28000
  'S25' = 1556
  'D26' = 44377
  'M26' = 30
  'N26' = 40000
  'O26' = 40000
  'P26' = 40000
  'S26' = 2222
  'D27' = 44432
  'M27' = 10
  'D28' = 44662
  'M28' = 15
  'N28' = 30000
  'O28' = 30000
  'P28' = 30000
  'S28' = 1667
  'D29' = 44449
  'N29' = 38000
  'O29' = 38000
  'P29' = 38000
  'S29' = 2111
  'D30' = 44379
  'M30' = 10
  'R30' = 'Región de Arica y Parinacota'
  'D31' = 44720
  'R31' = 'Perú'
  'D32' = 44294
  'M32' = 15
  'R32' = 'Región de Arica y Parinacota'
  'M33' = 10
  'R33' = 'Perú'
  'D34' = 44435
  'M34' = 105
  'D35' = 44357
  'M35' = 10
  'N35' = 38000
  'O35' = 38000
  'P35' = 38000
  'R35' = 'Perú'
  'S35' = 2111
  'D37' = 44726
  'M37' = 30
  'N37' = 34000
  'O37' = 34000
  'P37' = 34000
  'S37' = 1889
  'D38' = 44699
  'M38' = 20
  'N38' = 35000
  'O38' = 35000
  'P38' = 35000
  'S38' = 1944
  'D39' = 44658
  'M39' = 30
  'N39' = 28000
  'O39' = 28000
  'P39' = 28000
  'S39' = 1556
  'D40' = 44690
  'M40' = 25
  'N40' = 34000
  'O40' = 34000
  'P40' = 34000
  'R40' = 'Región de Arica y Parinacota'
  'S40' = 1889
  'D41' = 44433
  'M41' = 15
  'N41' = 35000
  'O41' = 35000
  'P41' = 35000
  'Q41' = '$/caja 18 kilos'
  'S41' = 1944
  'T41' = 18
  'D42' = 44442
  'M42' = 15
  'N42' = 35000
  'O42' = 35000
  'P42' = 35000
  'R42' = 'Perú'
  'S42' = 1944
}

foreach ($key in $updates.Keys) {
  $ws.Range($key).Value = $updates[$key]
}
